$d = $word.ActiveDocument

# Locate the target paragraph - the "Θέμα: «Έγκριση μετακίνησης ..." heading,
# uniquely identified by the word "σχολείου" which only appears here.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains('σχολείου')) {
        $target = $p
    }
}

# Step 1: "... εκπαιδευτικών του σχολείου: «${" -> "... εκπαιδευτικών του: «${"
$r1 = $target.Range
$r1.Find.Execute('εκπαιδευτικών του σχολείου: «${', $true, $false, $false, $false, $false, $true, 1, $false, 'εκπαιδευτικών του: «${', 2)

# Step 2: wrap the first "${country}" placeholder (the one in this paragraph)
# with guillemets: " ${country}," -> " «${country}»,"
$r2 = $target.Range
$r2.Find.Execute('${country}', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$parent = $r2.Find.Parent
$insAfter = $d.Range($parent.End, $parent.End)
$insAfter.InsertBefore([char]0xBB)
$insBefore = $d.Range($parent.Start, $parent.Start)
$insBefore.InsertBefore([char]0xAB)

# Step 3: delete the trailing red, size-22(half-pt)/11pt space run at the very
# end of the paragraph (just before the paragraph mark).
$r3 = $target.Range
$endPos = $r3.End
$lastCharRange = $d.Range($endPos - 2, $endPos - 1)
if ($lastCharRange.Text -eq " " -and $lastCharRange.Font.Color -eq 255) {
    $lastCharRange.Delete()
}
